$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A174").Value = "29.03.2024"
$ws.Range("B174").Value = "29.03.2024"
$ws.Range("C174").Value = 32623
$ws.Range("D174").Value = 13000
$ws.Range("E174").Value = 8400
$ws.Range("F174").Value = 75092
$ws.Range("G174").Value = 8663
$ws.Range("H174").Value = 6327
$ws.Range("I174").Value = 8000
$ws.Range("J174").Value = 454
$ws.Range("K174").Value = 116
$ws.Range("L174").Value = 4750
$ws.Range("M174").Value = "https://web.archive.org/web/20240329155650/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

$ws.Range("A175").Value = "30.03.2024"
$ws.Range("B175").Value = "29.03.2024"
$ws.Range("C175").Value = 32623
$ws.Range("D175").Value = 13000
$ws.Range("E175").Value = 8400
$ws.Range("F175").Value = 75092
$ws.Range("G175").Value = 8663
$ws.Range("H175").Value = 6327
$ws.Range("I175").Value = 8000
$ws.Range("J175").Value = 454
$ws.Range("K175").Value = 116
$ws.Range("L175").Value = 4750
$ws.Range("M175").Value = "https://web.archive.org/web/20240330182428/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

$ws.Range("A176").Value = "31.03.2024"
$ws.Range("B176").Value = "31.03.2024"
$ws.Range("C176").Value = 32782
$ws.Range("D176").Value = 13000
$ws.Range("E176").Value = 8400
$ws.Range("F176").Value = 75298
$ws.Range("G176").Value = 8663
$ws.Range("H176").Value = 6327
$ws.Range("I176").Value = 8000
$ws.Range("J176").Value = 455
$ws.Range("K176").Value = 116
$ws.Range("L176").Value = 4750
$ws.Range("M176").Value = "https://web.archive.org/web/20240331170708/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"
